# Weekly refresh of Fruta/Hortaliza price data (Granada, Vega Central Mapocho de Santiago).
# Columns D, K, L, M, N, O, P, Q, R, S, T are refreshed per row; A-C, E-J (Mercado/Producto
# identification columns) are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, Fecha(D), Variedad(K), Calidad(L), Volumen(M), PrecioMin(N), PrecioMax(O), PrecioProm(P), Unidad(Q), Origen(R), PrecioKg(S), KgUnidad(T)
$rows = @(
    ,@(2, 44285, 'Wonderfull', 'Especial', 40, 18000, 18000, 18000, '$/caja 15 kilos empedrada', 'Provincia del Elquí', 1200, 15)
    ,@(3, 44285, 'Wonderfull', 'Primera', 90, 15000, 15000, 15000, '$/caja 15 kilos empedrada', 'Provincia del Elquí', 1000, 15)
    ,@(4, 44285, 'Wonderfull', 'Segunda', 75, 12000, 12000, 12000, '$/caja 15 kilos empedrada', 'Provincia del Elquí', 800, 15)
    ,@(5, 44687, 'Wonderfull', 'Especial', 220, 21000, 21000, 21000, '$/caja 18 kilos granel', 'Región de O''Higgins', 1167, 18)
    ,@(6, 44687, 'Wonderfull', 'Primera', 250, 15000, 15000, 15000, '$/caja 18 kilos granel', 'Región de O''Higgins', 833, 18)
    ,@(7, 44687, 'Wonderfull', 'Segunda', 280, 10000, 10000, 10000, '$/caja 18 kilos granel', 'Región de O''Higgins', 556, 18)
    ,@(8, 44694, 'Wonderfull', 'Especial', 200, 21600, 21600, 21600, '$/caja 18 kilos granel', 'Región de O''Higgins', 1200, 18)
    ,@(9, 44694, 'Wonderfull', 'Primera', 220, 18000, 18000, 18000, '$/caja 18 kilos granel', 'Región de O''Higgins', 1000, 18)
    ,@(10, 44694, 'Wonderfull', 'Segunda', 250, 14400, 14400, 14400, '$/caja 18 kilos granel', 'Región de O''Higgins', 800, 18)
    ,@(11, 44305, 'Wonderfull', 'Primera', 50, 18000, 18000, 18000, '$/caja 15 kilos granel', 'Región de O''Higgins', 1200, 15)
    ,@(12, 44305, 'Wonderfull', 'Segunda', 60, 15000, 15000, 15000, '$/caja 15 kilos granel', 'Región de O''Higgins', 1000, 15)
    ,@(13, 44698, 'Wonderfull', 'Especial', 280, 18000, 18000, 18000, '$/caja 15 kilos granel', 'Región de O''Higgins', 1200, 15)
    ,@(14, 44698, 'Wonderfull', 'Primera', 220, 15000, 15000, 15000, '$/caja 15 kilos granel', 'Región de O''Higgins', 1000, 15)
    ,@(15, 44698, 'Wonderfull', 'Segunda', 200, 12000, 12000, 12000, '$/caja 15 kilos granel', 'Región de O''Higgins', 800, 15)
    ,@(16, 44649, 'Sin especificar', 'Especial', 220, 21600, 21600, 21600, '$/caja 18 kilos granel', 'Provincia de Limarí', 1200, 18)
    ,@(17, 44649, 'Sin especificar', 'Primera', 250, 16200, 16200, 16200, '$/caja 18 kilos granel', 'Provincia de Limarí', 900, 18)
    ,@(18, 44649, 'Sin especificar', 'Segunda', 180, 14400, 14400, 14400, '$/caja 18 kilos granel', 'Provincia de Limarí', 800, 18)
    ,@(19, 44658, 'Sin especificar', 'Especial', 280, 21600, 21600, 21600, '$/caja 18 kilos granel', 'Provincia de Limarí', 1200, 18)
    ,@(20, 44658, 'Sin especificar', 'Primera', 330, 16200, 16200, 16200, '$/caja 18 kilos granel', 'Provincia de Limarí', 900, 18)
    ,@(21, 44658, 'Sin especificar', 'Segunda', 220, 14400, 14400, 14400, '$/caja 18 kilos granel', 'Provincia de Limarí', 800, 18)
    ,@(22, 44644, 'Sin especificar', 'Especial', 180, 18000, 18000, 18000, '$/caja 15 kilos granel', 'Provincia de Limarí', 1200, 15)
    ,@(23, 44644, 'Sin especificar', 'Primera', 220, 13500, 13500, 13500, '$/caja 15 kilos granel', 'Provincia de Limarí', 900, 15)
    ,@(24, 44644, 'Sin especificar', 'Segunda', 290, 12000, 12000, 12000, '$/caja 15 kilos granel', 'Provincia de Limarí', 800, 15)
    ,@(25, 44706, 'Wonderfull', 'Especial', 200, 16000, 16000, 16000, '$/caja 18 kilos granel', 'Región de O''Higgins', 889, 18)
    ,@(26, 44706, 'Wonderfull', 'Primera', 220, 12500, 12500, 12500, '$/caja 18 kilos granel', 'Región de O''Higgins', 694, 18)
    ,@(27, 44664, 'Sin especificar', 'Especial', 300, 21600, 21600, 21600, '$/caja 18 kilos granel', 'Provincia de Limarí', 1200, 18)
    ,@(28, 44664, 'Sin especificar', 'Primera', 250, 18000, 18000, 18000, '$/caja 18 kilos granel', 'Provincia de Limarí', 1000, 18)
    ,@(29, 44664, 'Sin especificar', 'Segunda', 250, 16000, 16000, 16000, '$/caja 18 kilos granel', 'Provincia de Limarí', 889, 18)
    ,@(30, 44678, 'Sin especificar', 'Especial', 290, 15000, 15000, 15000, '$/caja 15 kilos granel', 'Región de O''Higgins', 1000, 15)
    ,@(31, 44678, 'Sin especificar', 'Primera', 220, 12000, 12000, 12000, '$/caja 15 kilos granel', 'Región de O''Higgins', 800, 15)
    ,@(32, 44685, 'Wonderfull', 'Especial', 350, 21000, 21000, 21000, '$/caja 18 kilos granel', 'Región de O''Higgins', 1167, 18)
    ,@(33, 44685, 'Wonderfull', 'Primera', 330, 15000, 15000, 15000, '$/caja 18 kilos granel', 'Región de O''Higgins', 833, 18)
    ,@(34, 44685, 'Wonderfull', 'Segunda', 280, 10000, 10000, 10000, '$/caja 18 kilos granel', 'Región de O''Higgins', 556, 18)
    ,@(35, 44309, 'Wonderfull', 'Primera', 40, 18000, 18000, 18000, '$/caja 15 kilos granel', 'Región de O''Higgins', 1200, 15)
    ,@(36, 44309, 'Wonderfull', 'Segunda', 70, 15000, 15000, 15000, '$/caja 15 kilos granel', 'Región de O''Higgins', 1000, 15)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value  = $row[1]   # D Fecha
    $ws.Cells.Item($r, 11).Value = $row[2]   # K Variedad
    $ws.Cells.Item($r, 12).Value = $row[3]   # L Calidad
    $ws.Cells.Item($r, 13).Value = $row[4]   # M Volumen
    $ws.Cells.Item($r, 14).Value = $row[5]   # N Precio minimo
    $ws.Cells.Item($r, 15).Value = $row[6]   # O Precio maximo
    $ws.Cells.Item($r, 16).Value = $row[7]   # P Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $row[8]   # Q Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $row[9]   # R Origen
    $ws.Cells.Item($r, 19).Value = $row[10]  # S Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $row[11]  # T Kg / unidad
}
